$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(20, 1).Value = "228  Game Provider Name = QuickSpin   Game Name =  Eastern Emeralds  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/EasternEmeralds.png"
$ws.Cells.Item(21, 1).Value = "245  Game Provider Name = QuickSpin   Game Name =  Fairy Gate  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/FairyGate.png"
$ws.Cells.Item(22, 1).Value = "258  Game Provider Name = TomHorn   Game Name =  Feng Fu  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/FengFu.png"
$ws.Cells.Item(23, 1).Value = "265  Game Provider Name = TomHorn   Game Name =  Fire 'n' Hot  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/FirenHot.png"
$ws.Cells.Item(24, 1).Value = "277  Game Provider Name = TomHorn   Game Name =  Flaming Fruit  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/FlamingFruit.png"
$ws.Cells.Item(25, 1).Value = "291  Game Provider Name = TomHorn   Game Name =  Frozen Queen  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/FrozenQueen.png"
$ws.Cells.Item(26, 1).Value = "311  Game Provider Name = TomHorn   Game Name =  Geisha's Fan  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/GeishasFun.png"
$ws.Cells.Item(27, 1).Value = "317  Game Provider Name = QuickSpin   Game Name =  Genies Touch  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_geniestouch.jpg"
$ws.Cells.Item(28, 1).Value = "332  Game Provider Name = TomHorn   Game Name =  Gold x  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/gold-x.png"
$ws.Cells.Item(29, 1).Value = "345  Game Provider Name = QuickSpin   Game Name =  Goldilocks & Wild Bears  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_goldilocks2.jpg"
$ws.Cells.Item(30, 1).Value = "346  Game Provider Name = QuickSpin   Game Name =  Gold Lab  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_goldlab.jpg"
$ws.Cells.Item(31, 1).Value = "376  Game Provider Name = TomHorn   Game Name =  Hot Blizzard  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/hot-blizzard.png"
$ws.Cells.Item(32, 1).Value = "381  Game Provider Name = QuickSpin   Game Name =  Hot Sync  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/HotSync.png"
$ws.Cells.Item(33, 1).Value = "385  Game Provider Name = TomHorn   Game Name =  Hot'n'Fruity  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/hotnfruity.png"
$ws.Cells.Item(34, 1).Value = "402  Game Provider Name = QuickSpin   Game Name =  Illuminous  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_illuminous.jpg"
$ws.Cells.Item(35, 1).Value = "407  Game Provider Name = TomHorn   Game Name =  Inca's Treasure  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/incas-treasure.png"
$ws.Cells.Item(36, 1).Value = "413  Game Provider Name = RelaxGaming   Game Name =  Iron Bank  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/ironbank_rg.jpg"
$ws.Cells.Item(37, 1).Value = "436  Game Provider Name = TomHorn   Game Name =  Joker Reelz  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/joker-reelz.png"
$ws.Cells.Item(38, 1).Value = "473  Game Provider Name = QuickSpin   Game Name =  Leprechaun Hills  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_leprechaunhills.jpg"
$ws.Cells.Item(39, 1).Value = "514  Game Provider Name = QuickSpin   Game Name =  Mayana  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_mayana.jpg"
$ws.Cells.Item(40, 1).Value = "527  Game Provider Name = QuickSpin   Game Name =  Mighty Arthur  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_mightyarthur.jpg"
$ws.Cells.Item(41, 1).Value = "533  Game Provider Name = RelaxGaming   Game Name =  Money Train 2  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/moneytrain2_rg.jpg"
$ws.Cells.Item(42, 1).Value = "534  Game Provider Name = TomHorn   Game Name =  Monkey 27  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/Monkey27.png"
$ws.Cells.Item(43, 1).Value = "538  Game Provider Name = TomHorn   Game Name =  Monster Madness  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/MonsterMadness.png"
$ws.Cells.Item(44, 1).Value = "544  Game Provider Name = QuickSpin   Game Name =  Mountain King  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/MountainKing.png"
$ws.Cells.Item(45, 1).Value = "560  Game Provider Name = QuickSpin   Game Name =  Northern Sky  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/NorthernSky.png"
$ws.Cells.Item(46, 1).Value = "571  Game Provider Name = TomHorn   Game Name =  Panda's Run  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/PandasRun.png"
$ws.Cells.Item(47, 1).Value = "584  Game Provider Name = QuickSpin   Game Name =  Phoenix Sun  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_phoenixsun.jpg"
$ws.Cells.Item(48, 1).Value = "585  Game Provider Name = QuickSpin   Game Name =  Pied Piper  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/PiedPiper.png"
$ws.Cells.Item(49, 1).Value = "593  Game Provider Name = QuickSpin   Game Name =  Pirates Charm  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/PiratesCharm.png"
$ws.Cells.Item(50, 1).Value = "625  Game Provider Name = QuickSpin   Game Name =  Razortooth  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_sabretooth.jpg"
$ws.Cells.Item(51, 1).Value = "628  Game Provider Name = TomHorn   Game Name =  Red Lights  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/RedLights.png"
$ws.Cells.Item(52, 1).Value = "679  Game Provider Name = TomHorn   Game Name =  Savannah King  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/SavannahKing.png"
$ws.Cells.Item(53, 1).Value = "680  Game Provider Name = TomHorn   Game Name =  Scratch Card  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/ScratchCard.png"
$ws.Cells.Item(54, 1).Value = "685  Game Provider Name = QuickSpin   Game Name =  Second Strike  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_secondstrike.jpg"
$ws.Cells.Item(55, 1).Value = "692  Game Provider Name = QuickSpin   Game Name =  Sevens High  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_sevens.jpg"
$ws.Cells.Item(56, 1).Value = "695  Game Provider Name = TomHorn   Game Name =  Shaolin's Tiger  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/ShaolinsTiger.jpg"
$ws.Cells.Item(57, 1).Value = "697  Game Provider Name = TomHorn   Game Name =  Sherlock in Bohemia  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/Sherlock.png"
$ws.Cells.Item(58, 1).Value = "701  Game Provider Name = TomHorn   Game Name =  Sizable Win  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/SizableWin.png"
$ws.Cells.Item(59, 1).Value = "705  Game Provider Name = TomHorn   Game Name =  Sky Barons  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/SkyBarons.png"
$ws.Cells.Item(60, 1).Value = "708  Game Provider Name = RelaxGaming   Game Name =  Snake Arena  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/snakearena_rg.jpg"
$ws.Cells.Item(61, 1).Value = "721  Game Provider Name = TomHorn   Game Name =  Spinball  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/spinball.png"
$ws.Cells.Item(62, 1).Value = "722  Game Provider Name = QuickSpin   Game Name =  Spinions Beach Party  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_spinions.jpg"
$ws.Cells.Item(63, 1).Value = "732  Game Provider Name = QuickSpin   Game Name =  Sticky Bandits  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/qso_stickybandits.jpg"
$ws.Cells.Item(64, 1).Value = "749  Game Provider Name = TomHorn   Game Name =  Sweet Crush  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/sweet-crush.png"
$ws.Cells.Item(65, 1).Value = "757  Game Provider Name = QuickSpin   Game Name =  Tales of Doctor Dolittle  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/TalesofDoctorDolittle.png"
$ws.Cells.Item(66, 1).Value = "762  Game Provider Name = TomHorn   Game Name =  The Cup  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/TheCup.png"
$ws.Cells.Item(67, 1).Value = "786  Game Provider Name = TomHorn   Game Name =  The Secret of BA  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/the-secret-of-ba.png"
$ws.Cells.Item(68, 1).Value = "794  Game Provider Name = TomHorn   Game Name =  Thrones Of Persia  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/ThronesOfPersia.png"
$ws.Cells.Item(69, 1).Value = "839  Game Provider Name = QuickSpin   Game Name =  Volcano Riches  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/VolcanoRiches.png"
$ws.Cells.Item(70, 1).Value = "850  Game Provider Name = RelaxGaming   Game Name =  Wild Chapo  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/wildchapo_rg.jpg"
$ws.Cells.Item(71, 1).Value = "852  Game Provider Name = QuickSpin   Game Name =  Wild Chase  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/WildChase.png"
$ws.Cells.Item(72, 1).Value = "866  Game Provider Name = TomHorn   Game Name =  Wild Weather  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/WildWeather.png"
$ws.Cells.Item(73, 1).Value = "874  Game Provider Name = QuickSpin   Game Name =  Wins of Fortune  cod = 404   src = https://resources.bet2win.vip/products/outcomebet/web/WinsofFortune.png"
$ws.Cells.Item(74, 1).Value = "878  Game Provider Name = TomHorn   Game Name =  Wolf Sierra  cod = 404   src = https://resources.bet2win.vip/products/tomhorn/web/wolf-sierra.png"
